$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hyperlinked emails - set the display text first, then attach the mailto
# hyperlink so the existing cell text is reused as the display text.
$ws.Range("K7").Value = "sampadakadam98@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K7"), "mailto:sampadakadam98@gmail.com")

$ws.Range("K9").Value = "bisenpooja5413@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K9"), "mailto:bisenpooja5413@gmail.com")

$ws.Range("K6").Value = "kumbharomkar45@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K6"), "mailto:kumbharomkar45@gmail.com")

$ws.Range("K8").Value = "pritirpatil187@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K8"), "mailto:pritirpatil187@gmail.com")

$ws.Range("K11").Value = "mahadikamit086@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K11"), "mailto:mahadikamit086@gmail.com")

# Plain text email (no hyperlink) for row 10 / Shivani Bhosale
$ws.Range("K10").Value = "shivanisb234@gmail.com"
